# Weekly refresh of the Espinaca (Terminal La Palmera de La Serena) series:
# a new week's observation is inserted at row 20, every subsequent
# observation shifts down one row, and the oldest observation (formerly
# row 157) is preserved by appending it as the new last row (158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the last row (157) before it gets overwritten by the shift ---
$lastD = $ws.Cells.Item(157, 4).Value2
$lastJ = $ws.Cells.Item(157, 10).Value2
$lastK = $ws.Cells.Item(157, 11).Value2
$lastM = $ws.Cells.Item(157, 13).Value2
$lastP = $ws.Cells.Item(157, 16).Value2

# --- shift rows 21..157 down from rows 20..156 (walk bottom-up so the
#     source cell of each copy hasn't been touched yet) ---
for ($r = 157; $r -ge 21; $r--) {
    $prev = $r - 1
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($prev, 4).Value2
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($prev, 10).Value2
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($prev, 11).Value2
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($prev, 13).Value2
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($prev, 16).Value2
}

# --- row 20 becomes the new, most recent observation ---
$ws.Cells.Item(20, 4).Value2 = 44473
$ws.Cells.Item(20, 10).Value2 = 2000
$ws.Cells.Item(20, 11).Value2 = 400
$ws.Cells.Item(20, 13).Value2 = 450
$ws.Cells.Item(20, 16).Value2 = 900

# --- append a new last row (158) carrying what used to be row 157 ---
$ws.Cells.Item(158, 1).Value2 = $ws.Cells.Item(157, 1).Value2
$ws.Cells.Item(158, 2).Value2 = $ws.Cells.Item(157, 2).Value2
$ws.Cells.Item(158, 3).Value2 = $ws.Cells.Item(157, 3).Value2
$ws.Cells.Item(158, 4).Value2 = $lastD
$ws.Cells.Item(158, 4).NumberFormat = $ws.Cells.Item(157, 4).NumberFormat
$ws.Cells.Item(158, 5).Value2 = $ws.Cells.Item(157, 5).Value2
$ws.Cells.Item(158, 6).Value2 = $ws.Cells.Item(157, 6).Value2
$ws.Cells.Item(158, 7).Value2 = $ws.Cells.Item(157, 7).Value2
$ws.Cells.Item(158, 8).Value2 = $ws.Cells.Item(157, 8).Value2
$ws.Cells.Item(158, 9).Value2 = $ws.Cells.Item(157, 9).Value2
$ws.Cells.Item(158, 10).Value2 = $lastJ
$ws.Cells.Item(158, 11).Value2 = $lastK
$ws.Cells.Item(158, 12).Value2 = $ws.Cells.Item(157, 12).Value2
$ws.Cells.Item(158, 13).Value2 = $lastM
$ws.Cells.Item(158, 14).Value2 = $ws.Cells.Item(157, 14).Value2
$ws.Cells.Item(158, 15).Value2 = $ws.Cells.Item(157, 15).Value2
$ws.Cells.Item(158, 16).Value2 = $lastP
$ws.Cells.Item(158, 17).Value2 = $ws.Cells.Item(157, 17).Value2
$ws.Cells.Item(158, 18).Value2 = $ws.Cells.Item(157, 18).Value2
